$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.450.36'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.493.37'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.510'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.348'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.33%  '
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.952.09'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '69.341.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '24.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.505.95'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '352.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.27%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").Value = '  -4.37%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.620.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.63'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0871'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.55'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +141.57%  '
$ws.Range("E32").Value = '  -2.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '440.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.42%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("B35").Value = 'PancakeSwap'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.71'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '154.21'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.95%  '
$ws.Range("E37").Value = '  -2.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.06'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.59%  '
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("E41").Value = '  -1.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.60'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.88%  '
$ws.Range("E43").Value = '  -1.60%  '
$ws.Range("E44").Value = '  -1.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '139.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.505'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0724'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.93%  '
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0924'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.50%  '
